$d = $word.ActiveDocument

# Step 1: insert new Author paragraphs
$pAuthor = $d.Paragraphs(3)
$pAuthor.Range.InsertParagraphAfter()
$a1 = $d.Paragraphs(4)
$a1.Range.Text = "Kapil Dev Soni"
$a1.Range.InsertParagraphAfter()
$a2 = $d.Paragraphs(5)
$a2.Range.Text = "Nobhojit Roy"
$a2.Range.InsertParagraphAfter()
$a3 = $d.Paragraphs(6)
$a3.Range.Text = "Martin Gerdin Wärnberg"

# Step 2: update FirstParagraph text
$rngFP = $d.Paragraphs(7).Range
$ok = $rngFP.Find.Execute(". After adjusting for injury severity, the paper reported higher risk of mortality among patients transferred by police vehicles. India, with a large trauma burden, has limited availability of ambulances consequently, other forms of transport, such as police vans are used", $true, $false, $false, $false, $false, $true, 1, $false, ". After adjusting for injury severity, the paper reported higher risk of mortality among patients transferred by police vehicles. India, like Malawi, is a low- and middle-income country and has a large burden of trauma", 2)
if (-not $ok) { throw "Replace failed: . After adjusting for injury severity, the paper r" }
$rngFP2 = $d.Paragraphs(7).Range
$ok = $rngFP2.Find.Execute("[2].", $true, $false, $false, $false, $false, $true, 1, $false, "[2]. Similarly, India has limited pre-hospital services such and different forms of transport, such as police vans are used to take trauma patients to the hospital [3, 4].", 2)
if (-not $ok) { throw "Replace failed: [2]." }

# Step 3: update BodyText intro paragraph
$rngBI = $d.Paragraphs(8).Range
$ok = $rngBI.Find.Execute("We used the TITCO-India data set, based on data from four tertiary care hospitals across urban India", $true, $false, $false, $false, $false, $true, 1, $false, "We used the Towards Improved Trauma Care Outcomes (TITCO) in India cohort, based on data from four tertiary care hospitals across urban India", 2)
if (-not $ok) { throw "Replace failed: We used the TITCO-India data set, based on data fr" }

# Step 4: insert new Table 1 / Table 2 block (Compact style) before the old block
$anchor = $d.Paragraphs(9)
$anchor.Range.InsertBefore("Table 1:Characterstics of directly admitted Vehicular Trauma in TITCO-India Data set and Purceell et al. 2020`rVariable`rTITCO_India`rPurcell_et_all_2020`rGender (female %)`r13.69`r19.8`rAge, years, mean (SD)`r35.4 (12.9)`r31.6 (15.9)`rInjury Severity, median (IQR)`r10 (8-14)`r8 (5–9)`rMode of Transport (%)`rAmbulance`r23.77`r–`rPolice Vehicle`r39.59`r9.4`rPrivate vehicle`r21.87`r68.8`rMotor Rickshaw, Taxi car`r14.52`r–`rMinibus`r–`r17.4`rOthers`r–`r4.4`rTime of Presentation, median (IQR)`r1 (0.6-2)`r1 (0-3)`rMortality (%)`r15.25`r7.8`r* Injury Severity Score`r† Malawi Trauma Score`rTable 2: Relative Risk of Mortality of patients brought by Police Vehicles`rX_`rTITCO_India_Dataset`rPurcell_et_al_2020`rRelative Risk`r1.529`r1.56`rCI at 95%`r(1.032, 2.321)`r(1.13–2.17)`rp-value`r0.03`r0.08`r* Adjusted for time to presentation, injury severity, and injury mechanism, Reference group Ambulance`r† Adjusted for time for presentation, Age, Sex, and Injury Severity`r‡ Refernce group Ambulance Reference group private vehicles`r")

# Step 5: delete the old Table 1 / Table 2 block
$oldStart = $d.Paragraphs(65)
$oldEnd = $d.Paragraphs(96)
$rngOld = $d.Range($oldStart.Range.Start, $oldEnd.Range.End)
$rngOld.Delete()

# Step 6: update the three closing BodyText paragraphs
$rngB1 = $d.Paragraphs(65).Range
$ok = $rngB1.Find.Execute("Of the 16000 patients in the TITCO-India data set, 1668 were adult vehicular trauma patients who were directly admitted to the study sites. of these complete data set was available for 1157 patients.A brief demographic profile the study cohort is given in Table 1. The overall mortality was 15.21 per cent. The most common mode of transport used for direct vehicular trauma patients was police vehicles (39.59%), followed by ambulance (23.77%), private vehicles (21.87%), and public transport such as motor rickshaws and cars (14.52%).", $true, $false, $false, $false, $false, $true, 1, $false, "Of the 16000 patients in the TITCO cohort, 1668 were adult vehicular trauma patients who were directly admitted to the study sites. Of these complete data was available for 1109 patients. A brief comparison of the patient profile of the Indian and the Malawi cohorts is given in Table 1. The overall mortality was 15.19 per cent. The most common mode of transport used for direct vehicular trauma patients was police vehicles (39.96%), followed by ambulance (23.51%), private vehicles (22.33%), and public transport such as motor rickshaws and cars (14.2%).", 2)
if (-not $ok) { throw "Replace failed: Of the 16000 patients in the TITCO-India data set," }
$rngB2 = $d.Paragraphs(66).Range
$ok = $rngB2.Find.Execute("The Poisson multivariate regression analysis showed that the relative risk of mortality for police vehicles was higher when compared to patients transported by ambulances (RR 1.5, 95% CI x–y, p = y), when adjusting for age, sex, and ISS. This was lower than the relative risk of mortality due to private vehicles or taxis and motor rickshaws when compared to ambulances (Table 2).", $true, $false, $false, $false, $false, $true, 1, $false, "The Poisson multivariate regression analysis showed that the relative risk of mortality for police vehicles was higher when compared to patients transported by ambulances (RR 1.52, 95% CI 1.03– 2.32, p value = 0.03), when adjusting for age, sex, and ISS. This was higher than the relative risk of mortality due to private vehicles or taxis and motor rickshaws when compared to ambulances (Table 2).", 2)
if (-not $ok) { throw "Replace failed: The Poisson multivariate regression analysis showe" }
$rngB3 = $d.Paragraphs(67).Range
$ok = $rngB3.Find.Execute("This analysis of the TITCO-India data set show similar results to the findings by the authors using data from Malawi. Therefore, there is a need to explore context-specific strategies such as training the police personnel to address the burden of trauma mortality in low-resource settings.", $true, $false, $false, $false, $false, $true, 1, $false, "Nearly one-fourth of the patients in the TITCO cohort arrived by ambulance, indicating the maturity of the pre-hospital care system existing in urban India. The relative risk of mortality of patients brought by police vehicles are very similar (1.53 vs 1.56) in both the countries. Our analysis underscores that the solutions suggested by the authors for Malawi may be be applied in the urban Indian setting. Therefore, there is a need to explore context-specific strategies such as training the police personnel to address the burden of trauma mortality in LMIC settings with underdeveloped pre-hospital system.", 2)
if (-not $ok) { throw "Replace failed: This analysis of the TITCO-India data set show sim" }

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
